$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "'VDVXC2QEJE"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "'1.4.2015 г. 00:00:00 ч."
$ws.Range("D2").Value = "'IrregularExpense"
$ws.Range("E2").Value = "'First Irregular Expense"

# Row 3
$ws.Range("A3").Value = "'QGZUYP5A74"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "'15.4.2015 г. 00:00:00 ч."
$ws.Range("D3").Value = "'RegularExpense"
$ws.Range("E3").Value = "'Second regular exp"

# Row 4
$ws.Range("A4").Value = "'MZ8LZS4DDI"
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = "'7.4.2015 г. 00:00:00 ч."
$ws.Range("D4").Value = "'IrregularExpense"
$ws.Range("E4").Value = "'Third"

# Row 5
$ws.Range("A5").Value = "'GMWOGER643"
$ws.Range("B5").Value = 666
$ws.Range("C5").Value = "'8.4.2015 г. 00:00:00 ч."
$ws.Range("D5").Value = "'IrregularExpense"
$ws.Range("E5").Value = "'Strategy"

# Row 6
$ws.Range("A6").Value = "'DW24LQH232"
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = "'7.4.2015 г. 00:00:00 ч."
$ws.Range("D6").Value = "'IrregularExpense"
$ws.Range("E6").Value = "'test excel"

# Row 7
$ws.Range("A7").Value = "'7OFMUXDQU9"
$ws.Range("B7").Value = 34
$ws.Range("C7").Value = "'7.4.2015 г. 00:00:00 ч."
$ws.Range("D7").Value = "'RegularExpense"
$ws.Range("E7").Value = "'2fsd"

# Row 8 (new row)
$ws.Range("A8").Value = "'QCUCJKC37W"
$ws.Range("B8").Value = 678
$ws.Range("C8").Value = "'10.4.2015 г. 00:00:00 ч."
$ws.Range("D8").Value = "'RegularExpense"
$ws.Range("E8").Value = "'Malko po-dylyg tekst"

# Row 9 (new row)
$ws.Range("A9").Value = "'ARKSWRKWDA"
$ws.Range("B9").Value = 900
$ws.Range("C9").Value = "'29.4.2015 г."
$ws.Range("D9").Value = "'RegularExpense"
$ws.Range("E9").Value = "'Muahahaha"

# Restore selection as in target (active cell C5, single-cell selection)
$ws.Range("C5").Select() | Out-Null
